$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing task row values to reflect new response codes
$ws.Range("A2").Value = "New Task Def 1"
$ws.Range("B2").Value = "NEWTD1"

# Add the new "tutorial_stream" header column
$ws.Range("S1").Value = "tutorial_stream"

# Reflect the new active selection
$ws.Range("S1").Select()
